$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly data refresh: insert 3 new "Ciruela" price rows at row 363 (pushing the
# existing rows 363-381 down to 366-384), mirroring the upstream CSV dump's
# newest week (2023-04-05) being prepended ahead of the older records.
$ws.Range("A363:A365").EntireRow.Insert()

$newRows = @(
    @{Row=363; Date=45021; Variedad="Angeleno"; Calidad="Especial"; Vol=18; Pmin=230000; Pmax=240000; Pprom=235000; Origen="Región de O'Higgins"; PrecioKg=522},
    @{Row=364; Date=45021; Variedad="Angeleno"; Calidad="Primera";  Vol=14; Pmin=200000; Pmax=210000; Pprom=205000; Origen="Región de O'Higgins"; PrecioKg=456},
    @{Row=365; Date=45021; Variedad="Angeleno"; Calidad="Segunda";  Vol=12; Pmin=170000; Pmax=180000; Pprom=175000; Origen="Región de O'Higgins"; PrecioKg=389}
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row,1).Value  = 8
    $ws.Cells.Item($row,2).Value  = "Terminal La Palmera de La Serena"
    $ws.Cells.Item($row,3).Value  = "Coquimbo"
    $ws.Cells.Item($row,4).Value  = $r.Date
    $ws.Cells.Item($row,5).Value  = 4
    $ws.Cells.Item($row,6).Value  = "Fruta"
    $ws.Cells.Item($row,7).Value  = 100103
    $ws.Cells.Item($row,8).Value  = "Frutos de hueso (carozo)"
    $ws.Cells.Item($row,9).Value  = 100103002
    $ws.Cells.Item($row,10).Value = "Ciruela"
    $ws.Cells.Item($row,11).Value = $r.Variedad
    $ws.Cells.Item($row,12).Value = $r.Calidad
    $ws.Cells.Item($row,13).Value = $r.Vol
    $ws.Cells.Item($row,14).Value = $r.Pmin
    $ws.Cells.Item($row,15).Value = $r.Pmax
    $ws.Cells.Item($row,16).Value = $r.Pprom
    $ws.Cells.Item($row,17).Value = "`$/bins (450 kilos)"
    $ws.Cells.Item($row,18).Value = $r.Origen
    $ws.Cells.Item($row,19).Value = $r.PrecioKg
    $ws.Cells.Item($row,20).Value = 450
}
